# FIX: Tipos de cambio CORREGIDOS con valores reales (12/Nov/2025)
#   TC Compra: 517.50 -> 494.00
#   TC Venta:  525.00 -> 508.00
# Updates formulas that convert CRC -> USD, the CONFIG sheet values/history,
# the notes on the EFECTIVO sheet, and the explanatory cell comments.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: CONFIG
# ---------------------------------------------------------------------------
$cfg = $wb.Worksheets.Item("CONFIG")

# Fecha Creación (record was re-saved at 00:28 instead of 00:24)
$cfg.Range("B5").Value = "13/11/2025 00:28"

# TC Compra / TC Venta (current values) + "Última Actualización" timestamps
$cfg.Range("B12").Value = 494
$cfg.Range("C12").Value = "13/11/2025 00:28"

$cfg.Range("B13").Value = 508
$cfg.Range("C13").Value = "13/11/2025 00:28"

# Historial de Tipos de Cambio - first (initial) row
$cfg.Range("B17").Value = 494
$cfg.Range("C17").Value = 508

# Updated comments explaining TC Compra / TC Venta (blank line between each
# paragraph, matching the original comment layout)
[void]$cfg.Range("B12").Comment.Text("💡 TIPO DE CAMBIO COMPRA`n`nCuántos colones TE DAN por `$1 USD`n`nEjemplo: Si banco te da ₡494.00 por `$1`n`n⚠️ Actualizar 1 vez/semana")
[void]$cfg.Range("B13").Comment.Text("💡 TIPO DE CAMBIO VENTA`n`nCuántos colones PAGAS por `$1 USD`n`nEjemplo: Si banco cobra ₡508.00 por `$1`n`n⚠️ Actualizar 1 vez/semana")

# ---------------------------------------------------------------------------
# Sheet: EFECTIVO
# ---------------------------------------------------------------------------
$efec = $wb.Worksheets.Item("EFECTIVO")

# TOTAL BANCOS (USD equivalente) - CRC -> USD conversion rate
$efec.Range("E14").Formula = "=SUMIF(D5:D13,""USD"",E5:E13)+SUMIF(D5:D13,""CRC"",E5:E13)/494.0"

# TARJETAS DE CRÉDITO - Equiv. USD Total per card
$efec.Range("E18").Formula = "=C18+(D18/494)"
$efec.Range("E19").Formula = "=C19+(D19/494)"
$efec.Range("E20").Formula = "=C20+(D20/494)"
$efec.Range("E21").Formula = "=C21+(D21/494)"
$efec.Range("E22").Formula = "=C22+(D22/494)"

# Notes at the bottom of the sheet
$efec.Range("A30").Value = "• Tipo de cambio usado para conversión: TC Compra ₡494.00 por `$1 USD (12/Nov/2025)"
$efec.Range("A31").Value = "• EFECTIVO NETO REAL: Bancos - Tarjetas (ver cálculo arriba)"
